# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (blank) column inserted
# before the existing "Late" column (column N), pushing the old
# N/O/P ("Late" / "Outstanding"-heading / "Outstanding") columns one
# slot to the right (O/P/Q). The sheet becomes the active tab/sheet of
# the workbook (it previously was "Transactions").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N->O, O->P, P->Q).
$ws.Range("N1").EntireColumn.Insert() | Out-Null

# The new column picks up a width matching column M (~10.71 chars).
$ws.Columns("N").ColumnWidth = 9.8

# Make "Repayment schedule" the active sheet/tab (was "Transactions"),
# with the selection left on S14.
$ws.Activate() | Out-Null
$ws.Range("S14").Select() | Out-Null
